{"js": "// Replace each old expression/date with its new counterpart.\n// Document order: title paragraph, then 25 populated table cells\n// (the remaining cells in the 5x10 grid are intentionally blank).\nconst replacements = [\n  [\"2025-05-10 Saturday\", \"2025-05-11 Sunday\"],\n  [\"200\u00f78=\", \"720\u00f76=\"],\n  [\"173\u00f78=\", \"627\u00f75=\"],\n  [\"775\u00f75=\", \"591\u00f74=\"],\n  [\"997\u00f75=\", \"897\u00f75=\"],\n  [\"509\u00f74=\", \"850\u00f75=\"],\n  [\"186\u00f73=\", \"720\u00f74=\"],\n  [\"417\u00f74=\", \"152\u00f79=\"],\n  [\"175\u00f77=\", \"606\u00f74=\"],\n  [\"239\u00f72=\", \"534\u00f79=\"],\n  [\"258\u00f77=\", \"658\u00f74=\"],\n  [\"728\u00f75=\", \"897\u00f79=\"],\n  [\"711\u00f73=\", \"460\u00f79=\"],\n  [\"491\u00f75=\", \"737\u00f74=\"],\n  [\"732\u00f73=\", \"617\u00f74=\"],\n  [\"116\u00f77=\", \"186\u00f72=\"],\n  [\"393\u00f77=\", \"566\u00f79=\"],\n  [\"261\u00f75=\", \"858\u00f77=\"],\n  [\"389\u00f79=\", \"465\u00f77=\"],\n  [\"288\u00f73=\", \"276\u00f75=\"],\n  [\"905\u00f76=\", \"279\u00f72=\"],\n  [\"154\u00f75=\", \"854\u00f77=\"],\n  [\"314\u00f76=\", \"719\u00f76=\"],\n  [\"845\u00f74=\", \"359\u00f75=\"],\n  [\"208\u00f78=\", \"863\u00f79=\"],\n  [\"638\u00f75=\", \"430\u00f76=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each old expression/date with its new counterpart.\n# Document order: title paragraph, then 25 populated table cells\n# (the remaining cells in the 5x10 grid are intentionally blank).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"2025-05-10 Saturday\"; New = \"2025-05-11 Sunday\" },\n    @{ Old = \"200\u00f78=\"; New = \"720\u00f76=\" },\n    @{ Old = \"173\u00f78=\"; New = \"627\u00f75=\" },\n    @{ Old = \"775\u00f75=\"; New = \"591\u00f74=\" },\n    @{ Old = \"997\u00f75=\"; New = \"897\u00f75=\" },\n    @{ Old = \"509\u00f74=\"; New = \"850\u00f75=\" },\n    @{ Old = \"186\u00f73=\"; New = \"720\u00f74=\" },\n    @{ Old = \"417\u00f74=\"; New = \"152\u00f79=\" },\n    @{ Old = \"175\u00f77=\"; New = \"606\u00f74=\" },\n    @{ Old = \"239\u00f72=\"; New = \"534\u00f79=\" },\n    @{ Old = \"258\u00f77=\"; New = \"658\u00f74=\" },\n    @{ Old = \"728\u00f75=\"; New = \"897\u00f79=\" },\n    @{ Old = \"711\u00f73=\"; New = \"460\u00f79=\" },\n    @{ Old = \"491\u00f75=\"; New = \"737\u00f74=\" },\n    @{ Old = \"732\u00f73=\"; New = \"617\u00f74=\" },\n    @{ Old = \"116\u00f77=\"; New = \"186\u00f72=\" },\n    @{ Old = \"393\u00f77=\"; New = \"566\u00f79=\" },\n    @{ Old = \"261\u00f75=\"; New = \"858\u00f77=\" },\n    @{ Old = \"389\u00f79=\"; New = \"465\u00f77=\" },\n    @{ Old = \"288\u00f73=\"; New = \"276\u00f75=\" },\n    @{ Old = \"905\u00f76=\"; New = \"279\u00f72=\" },\n    @{ Old = \"154\u00f75=\"; New = \"854\u00f77=\" },\n    @{ Old = \"314\u00f76=\"; New = \"719\u00f76=\" },\n    @{ Old = \"845\u00f74=\"; New = \"359\u00f75=\" },\n    @{ Old = \"208\u00f78=\"; New = \"863\u00f79=\" },\n    @{ Old = \"638\u00f75=\"; New = \"430\u00f76=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n    if (-not $result) {\n        throw \"Text not found: $($pair.Old)\"\n    }\n}\n"}
